$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# sheet1 ("files"): add a new "demand" source row (row 7), mirroring the
# existing "state_gen" row (row 6), with its own hyperlink.
# ---------------------------------------------------------------------------
$ws1.Cells.Item(7,1).Value = "demand"

# Add the hyperlink first (it stamps the cell text with the display string),
# then overwrite the cell's displayed value with the real "S:\DEMAND" text
# so the two stay independent, exactly like the other rows in this sheet.
$ws1.Hyperlinks.Add($ws1.Cells.Item(7,2), "file:///\\10.2.100.51\scada\Reports\State_gen", $null, $null, "\\10.2.100.51\scada\Reports\State_gen") | Out-Null
$ws1.Cells.Item(7,2).Value = "S:\DEMAND"
$ws1.Cells.Item(7,2).Style = "Hyperlink"

$ws1.Cells.Item(7,3).Value = "Demand_"
$ws1.Cells.Item(7,4).Value = "%d_%m_%Y"
$ws1.Cells.Item(7,4).Font.Color = $ws1.Cells.Item(6,4).Font.Color
$ws1.Cells.Item(7,5).Value = ".csv"
$ws1.Cells.Item(7,6).Value = "csv"
$ws1.Cells.Item(7,7).Value = 2
$ws1.Cells.Item(7,8).Value = 6
$ws1.Cells.Item(7,9).Value = 1
$ws1.Cells.Item(7,10).Value = "%d-%m-%Y %H:%M:%S"

# ---------------------------------------------------------------------------
# sheet2 ("file_meas_info"): append 40 new "demand" measurement rows
# (rows 860-899).
# ---------------------------------------------------------------------------
$ws2.Cells.Item(860,1).Value = "demand"
$ws2.Cells.Item(860,2).Value = 2
$ws2.Cells.Item(860,3).Value = "WRLDCMP.SCADA1.A0036324"
$ws2.Cells.Item(860,4).Value = "SUBSTN.SYSCA_WR.SYSTEM.RAW_FREQ.MEAS.HZ"
$ws2.Cells.Item(861,1).Value = "demand"
$ws2.Cells.Item(861,2).Value = 3
$ws2.Cells.Item(861,3).Value = "WRLDCMP.SCADA1.A0047000"
$ws2.Cells.Item(861,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.WR_TOT_DEMAND.ADD.MW"
$ws2.Cells.Item(862,1).Value = "demand"
$ws2.Cells.Item(862,2).Value = 4
$ws2.Cells.Item(862,3).Value = "WRLDCMP.SCADA1.A0046980"
$ws2.Cells.Item(862,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MSEB_DEMAND.ADD.MW"
$ws2.Cells.Item(863,1).Value = "demand"
$ws2.Cells.Item(863,2).Value = 5
$ws2.Cells.Item(863,3).Value = "WRLDCMP.SCADA3.A0101733"
$ws2.Cells.Item(863,4).Value = "RCTMH_WR.SYSTEM.MUM_DEMAND.MEAS.MW"
$ws2.Cells.Item(864,1).Value = "demand"
$ws2.Cells.Item(864,2).Value = 6
$ws2.Cells.Item(864,3).Value = "WRLDCMP.SCADA1.A0046957"
$ws2.Cells.Item(864,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.GEB_DEMAND.ADD.MW"
$ws2.Cells.Item(865,1).Value = "demand"
$ws2.Cells.Item(865,2).Value = 7
$ws2.Cells.Item(865,3).Value = "WRLDCMP.SCADA1.A0046978"
$ws2.Cells.Item(865,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MPSEB_DEMAND.ADD.MW"
$ws2.Cells.Item(866,1).Value = "demand"
$ws2.Cells.Item(866,2).Value = 8
$ws2.Cells.Item(866,3).Value = "WRLDCMP.SCADA1.A0046945"
$ws2.Cells.Item(866,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.CSEB_DEMAND.ADD.MW"
$ws2.Cells.Item(867,1).Value = "demand"
$ws2.Cells.Item(867,2).Value = 9
$ws2.Cells.Item(867,3).Value = "WRLDCMP.SCADA1.A0046962"
$ws2.Cells.Item(867,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.GOA_DEMAND.ADD.MW"
$ws2.Cells.Item(868,1).Value = "demand"
$ws2.Cells.Item(868,2).Value = 10
$ws2.Cells.Item(868,3).Value = "WRLDCMP.SCADA1.A0046948"
$ws2.Cells.Item(868,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.DD_DEMAND.ADD.MW"
$ws2.Cells.Item(869,1).Value = "demand"
$ws2.Cells.Item(869,2).Value = 11
$ws2.Cells.Item(869,3).Value = "WRLDCMP.SCADA1.A0046953"
$ws2.Cells.Item(869,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.DNH_DEMAND.ADD.MW"
$ws2.Cells.Item(870,1).Value = "demand"
$ws2.Cells.Item(870,2).Value = 12
$ws2.Cells.Item(870,3).Value = "WRLDCMP.SCADA1.A0046960"
$ws2.Cells.Item(870,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.GEB_TOT_GEN.ADD.MW"
$ws2.Cells.Item(871,1).Value = "demand"
$ws2.Cells.Item(871,2).Value = 13
$ws2.Cells.Item(871,3).Value = "WRLDCMP.SCADA1.A0046984"
$ws2.Cells.Item(871,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MSEB_TOT_GEN.ADD.MW"
$ws2.Cells.Item(872,1).Value = "demand"
$ws2.Cells.Item(872,2).Value = 14
$ws2.Cells.Item(872,3).Value = "WRLDCMP.SCADA1.A0046979"
$ws2.Cells.Item(872,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MPSEB_TOT_GEN.ADD.MW"
$ws2.Cells.Item(873,1).Value = "demand"
$ws2.Cells.Item(873,2).Value = 15
$ws2.Cells.Item(873,3).Value = "WRLDCMP.SCADA1.A0046947"
$ws2.Cells.Item(873,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.CSEB_TOT_GEN.ADD.MW"
$ws2.Cells.Item(874,1).Value = "demand"
$ws2.Cells.Item(874,2).Value = 16
$ws2.Cells.Item(874,3).Value = "WRLDCMP.SCADA1.A0047001"
$ws2.Cells.Item(874,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.WR_TOT_GEN.ADD.MW"
$ws2.Cells.Item(875,1).Value = "demand"
$ws2.Cells.Item(875,2).Value = 17
$ws2.Cells.Item(875,3).Value = "WRLDCMP.SCADA1.A0046958"
$ws2.Cells.Item(875,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.GEB_HYDRO_GEN.ADD.MW"
$ws2.Cells.Item(876,1).Value = "demand"
$ws2.Cells.Item(876,2).Value = 18
$ws2.Cells.Item(876,3).Value = "WRLDCMP.SCADA1.A0046959"
$ws2.Cells.Item(876,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.GEB_TH_GS_GEN.ADD.MW"
$ws2.Cells.Item(877,1).Value = "demand"
$ws2.Cells.Item(877,2).Value = 19
$ws2.Cells.Item(877,3).Value = "WRLDCMP.SCADA1.A0046981"
$ws2.Cells.Item(877,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MSEB_HYD_GEN.ADD.MW"
$ws2.Cells.Item(878,1).Value = "demand"
$ws2.Cells.Item(878,2).Value = 20
$ws2.Cells.Item(878,3).Value = "WRLDCMP.SCADA1.A0046982"
$ws2.Cells.Item(878,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MSEB_TH_GS_GEN.ADD.MW"
$ws2.Cells.Item(879,1).Value = "demand"
$ws2.Cells.Item(879,2).Value = 21
$ws2.Cells.Item(879,3).Value = "WRLDCMP.SCADA1.A0047287"
$ws2.Cells.Item(879,4).Value = "SUBSTN.MPRTC_MP.STTN.TOTHYD_MW.MEAS.MW"
$ws2.Cells.Item(880,1).Value = "demand"
$ws2.Cells.Item(880,2).Value = 22
$ws2.Cells.Item(880,3).Value = "WRLDCMP.SCADA1.A0046992"
$ws2.Cells.Item(880,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.THRML_MW.ADD.MW"
$ws2.Cells.Item(881,1).Value = "demand"
$ws2.Cells.Item(881,2).Value = 23
$ws2.Cells.Item(881,3).Value = "WRLDCMP.SCADA1.A0046946"
$ws2.Cells.Item(881,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.CSEB_HYDRO_MW.MEAS.MW"
$ws2.Cells.Item(882,1).Value = "demand"
$ws2.Cells.Item(882,2).Value = 24
$ws2.Cells.Item(882,3).Value = "WRLDCMP.SCADA1.A0046391"
$ws2.Cells.Item(882,4).Value = "SUBSTN.SYSCA_CG.SYSTEM.ALL_GEN_CG_MW.MEAS.MW"
$ws2.Cells.Item(883,1).Value = "demand"
$ws2.Cells.Item(883,2).Value = 25
$ws2.Cells.Item(883,3).Value = "WRLDCMP.SCADA1.A0046999"
$ws2.Cells.Item(883,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.WR_NUCLEAR_GEN.ADD.MW"
$ws2.Cells.Item(884,1).Value = "demand"
$ws2.Cells.Item(884,2).Value = 26
$ws2.Cells.Item(884,3).Value = "WRLDCMP.SCADA1.A0047002"
$ws2.Cells.Item(884,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.WR_TOT_HYD_GEN.ADD.MW"
$ws2.Cells.Item(885,1).Value = "demand"
$ws2.Cells.Item(885,2).Value = 27
$ws2.Cells.Item(885,3).Value = "WRLDCMP.SCADA1.A0047003"
$ws2.Cells.Item(885,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.WR_TOT_TH_GS.ADD.MW"
$ws2.Cells.Item(886,1).Value = "demand"
$ws2.Cells.Item(886,2).Value = 28
$ws2.Cells.Item(886,3).Value = "WRLDCMP.SCADA1.A0046961"
$ws2.Cells.Item(886,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.GEB_TRML_GEN.ADD.MW"
$ws2.Cells.Item(887,1).Value = "demand"
$ws2.Cells.Item(887,2).Value = 29
$ws2.Cells.Item(887,3).Value = "WRLDCMP.SCADA1.A0046983"
$ws2.Cells.Item(887,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.MSEB_THRM_GEN.ADD.MW"
$ws2.Cells.Item(888,1).Value = "demand"
$ws2.Cells.Item(888,2).Value = 30
$ws2.Cells.Item(888,3).Value = "WRLDCMP.SCADA1.A0047004"
$ws2.Cells.Item(888,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.WR_TOT_TH_ONLY.ADD.MW"
$ws2.Cells.Item(889,1).Value = "demand"
$ws2.Cells.Item(889,2).Value = 31
$ws2.Cells.Item(889,3).Value = "WRLDCMP.SCADA1.A0049839"
$ws2.Cells.Item(889,4).Value = "SUBSTN.ATHN4_IP.STTN.USER_CALC.MEAS.UC1"
$ws2.Cells.Item(890,1).Value = "demand"
$ws2.Cells.Item(890,2).Value = 32
$ws2.Cells.Item(890,3).Value = "WRLDCMP.SCADA3.A0108222"
$ws2.Cells.Item(890,4).Value = "SUBSTN.PGVCL_GJ.CALC.PG_UG_SLR_TTL.ADD.MW"
$ws2.Cells.Item(891,1).Value = "demand"
$ws2.Cells.Item(891,2).Value = 33
$ws2.Cells.Item(891,3).Value = "WRLDCMP.SCADA3.A0104731"
$ws2.Cells.Item(891,4).Value = "SUBSTN.PGVCL_GJ.CALC.PG_WIND_TTL.ADD.MW"
$ws2.Cells.Item(892,1).Value = "demand"
$ws2.Cells.Item(892,2).Value = 34
$ws2.Cells.Item(892,3).Value = "WRLDCMP.SCADA3.A0108546"
$ws2.Cells.Item(892,4).Value = "SUBSTN.MPRTC_MP.STTN.RE_SOL_TOT.MEAS.MW"
$ws2.Cells.Item(893,1).Value = "demand"
$ws2.Cells.Item(893,2).Value = 35
$ws2.Cells.Item(893,3).Value = "WRLDCMP.SCADA3.A0108547"
$ws2.Cells.Item(893,4).Value = "SUBSTN.MPRTC_MP.STTN.RE_WIND_TOT.MEAS.MW"
$ws2.Cells.Item(894,1).Value = "demand"
$ws2.Cells.Item(894,2).Value = 36
$ws2.Cells.Item(894,3).Value = "WRLDCMP.SCADA3.A0103074"
$ws2.Cells.Item(894,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.CG_SOLAR_TOT.ADD.MW"
$ws2.Cells.Item(895,1).Value = "demand"
$ws2.Cells.Item(895,2).Value = 37
$ws2.Cells.Item(895,3).Value = "WRLDCMP.SCADA3.A0109531"
$ws2.Cells.Item(895,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.TOT_WIND_MH.ADDT.MW"
$ws2.Cells.Item(896,1).Value = "demand"
$ws2.Cells.Item(896,2).Value = 38
$ws2.Cells.Item(896,3).Value = "WRLDCMP.SCADA3.A0108220"
$ws2.Cells.Item(896,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.TOT_SOLAR_MH.ADD.MW"
$ws2.Cells.Item(897,1).Value = "demand"
$ws2.Cells.Item(897,2).Value = 39
$ws2.Cells.Item(897,3).Value = "WRLDCMP.SCADA3.A0106536"
$ws2.Cells.Item(897,4).Value = "SUBSTN.IMEXP_WR.ABT.ALL_INDIA_DMD.MEAS.MW"
$ws2.Cells.Item(898,1).Value = "demand"
$ws2.Cells.Item(898,2).Value = 40
$ws2.Cells.Item(898,3).Value = "WRLDCMP.SCADA1.A0043296"
$ws2.Cells.Item(898,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.ESIL_TOT_DRWL.ADD.MW"
$ws2.Cells.Item(899,1).Value = "demand"
$ws2.Cells.Item(899,2).Value = 41
$ws2.Cells.Item(899,3).Value = "WRLDCMP.SCADA3.A0111629"
$ws2.Cells.Item(899,4).Value = "SUBSTN.IMEXP_WR.SYSTEM.BARC_TOTAL.ADDT.MW"

# ---------------------------------------------------------------------------
# View state: sheet2 ("file_meas_info") becomes the active/selected sheet,
# with its selection reset to C1 and sheet1's selection moved to C15.
# ---------------------------------------------------------------------------
$ws1.Range("C15").Select()
$ws2.Activate()
$ws2.Range("C1").Select()
